$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("E1").Value = "Delta NPV"
$ws.Range("F1").Value = "Std Delta NPV (upper bound)"

# Row 2 (baseline) - Delta NPV is 0 by definition, Std Delta upper bound is "-"
$ws.Range("E2").Formula = "=B2-`$B`$2"
$ws.Range("F2").Value = "-"

$ws.Range("G1").Value = "Std Delta NPV (%)"

# Rows 3-9: Delta NPV, Std Delta NPV (upper bound), Std Delta NPV (%)
for ($r = 3; $r -le 9; $r++) {
    $ws.Range("E$r").Formula = "=B$r-`$B`$2"
    $ws.Range("F$r").Formula = "=SQRT(POWER(`$C`$2,2)+POWER(C$r,2))"
    $ws.Range("G$r").Formula = "=F$r/E$r"
    $ws.Range("G$r").NumberFormat = "0.00%"
}

$ws.Range("I9").Select()
